# Edit HSI_Review.xlsx per commit "Update  HSI review sheet"
$wb = $excel.ActiveWorkbook

$wsIntro = $wb.Worksheets.Item("Introduction ")
$wsCross = $wb.Worksheets.Item("Cross review points ")

# --- Introduction sheet edits ---
# D7: version number 1 -> 1.5
$wsIntro.Range("D7").Value = 1.5

# D9: Last update date, was text "24/01/2020" -> real date 02/09/2020
$wsIntro.Range("D9").Value = 44076
$wsIntro.Range("D9").NumberFormat = "mm-dd-yy"

# New history row 14: version 0.2 entry
$wsIntro.Range("B14").Value = 0.2
$wsIntro.Range("C14").Value = "T.Sharaby"
$wsIntro.Range("E14").Value = 44076
$wsIntro.Range("E14").NumberFormat = "mm-dd-yy"
$wsIntro.Range("G14").Value = "Update the status for each req "

$wsIntro.Range("G15").Select()

# --- Cross review points sheet edits ---
$wsCross.Range("H2").Value = "Resolved"
$wsCross.Range("H3").Value = "Resolved"
$wsCross.Range("H4").Value = "Resolved"
$wsCross.Range("H5").Value = "Resolved"
$wsCross.Range("H6").Value = "Resolved"

$wsCross.Range("F15").Select()

# Make Introduction the active sheet (matches final tabSelected state)
$wsIntro.Select()
